$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Set column A width (target OOXML width 21.3368944440569; closest reachable value via
# the character-based ColumnWidth COM property, which this runtime quantizes to 1/6-pt
# pixel steps, is 21.333333333333332)
$ws.Columns.Item(1).ColumnWidth = 20.55

# Add new rows of data (rows 2-4)
$ws.Range("A2").Value = "Cumplimiento de pago"
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 3
$ws.Range("D2").Value = 4
$ws.Range("E2").Value = 8
$ws.Range("F2").Value = 7

$ws.Range("A3").Value = "Check-in"
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 4
$ws.Range("D3").Value = 6
$ws.Range("E3").Value = 5
$ws.Range("F3").Value = 4

$ws.Range("A4").Value = "Messi"
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = 4
$ws.Range("D4").Value = 4
$ws.Range("E4").Value = 4
$ws.Range("F4").Value = 6

# Update the chart series formulas (categories & values)
$chartObj = $ws.ChartObjects(1)
$chart = $chartObj.Chart

$valCols = @("B", "C", "D", "E", "F")
for ($i = 1; $i -le 5; $i++) {
    $series = $chart.SeriesCollection($i)
    $series.XValues = "Sheet1!A2:A4"
    $colLetter = $valCols[$i - 1]
    $series.Values = "Sheet1!" + $colLetter + "2:" + $colLetter + "4"
}

# Resize the chart (to.col from 7 to 6, to.colOff from 495300 to 295275)
$chartObj.Width = 397.69140625
